$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# The "id" row in the data dictionary used to carry a bold red-highlight
# label ("ID of the participant"). Re-label it to the shorter "ID " and
# drop the red/bold emphasis so the row matches the plain styling used
# for the rest of the dictionary (borrow that plain style from a cell
# that already carries it, e.g. Categories!C2, so we don't mint a new
# font/style record).
$wsCategories.Range("C2").Copy()
$ws.Range("B2:D2").PasteSpecial(-4122)
$ws.Range("C2").Value = "ID "

$excel.CutCopyMode = $false
